$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price, volume and hour values for cryptos list (Feb 13 2023 update)
$cells = @{
    "D2" = "312.93"
    "E2" = "1.82%"
    "G2" = "5"
    "D3" = "40.00"
    "E3" = "-2.57%"
    "G3" = "5"
    "D4" = "5.182"
    "E4" = "-1.35%"
    "G4" = "5"
    "D5" = "0.07581"
    "E5" = "-1.05%"
    "G5" = "5"
    "D6" = "4.329"
    "E6" = "0.14%"
    "G6" = "5"
    "D7" = "1.669"
    "E7" = "2.94%"
    "G7" = "5"
    "D8" = "0.9262"
    "E8" = "0.88%"
    "G8" = "5"
    "G9" = "5"
    "D10" = "0.1199"
    "E10" = "-4.39%"
    "G10" = "5"
    "D11" = "0.1817"
    "E11" = "-1.01%"
    "G11" = "5"
    "D12" = "0.09092"
    "E12" = "-0.92%"
    "G12" = "5"
    "D13" = "0.04132"
    "E13" = "-3.05%"
    "G13" = "5"
    "D14" = "0.1052"
    "E14" = "0.07%"
    "G14" = "5"
    "D15" = "0.001281"
    "E15" = "1.34%"
    "G15" = "5"
    "E16" = "1.53%"
    "G16" = "5"
    "G17" = "5"
    "D18" = "3.350"
    "E18" = "-0.08%"
    "G18" = "5"
    "D19" = "0.3357"
    "E19" = "0.63%"
    "G19" = "5"
    "D20" = "7.593"
    "E20" = "6.09%"
    "G20" = "5"
    "D21" = "0.1351"
    "E21" = "-2.40%"
    "G21" = "5"
    "E22" = "-2.89%"
    "G22" = "5"
    "D23" = "0.04016"
    "E23" = "-1.26%"
    "G23" = "5"
    "D24" = "0.001276"
    "E24" = "0.89%"
    "G24" = "5"
    "D25" = "0.003976"
    "E25" = "-4.53%"
    "G25" = "5"
    "E26" = "-0.30%"
    "G26" = "5"
    "G27" = "5"
    "G28" = "5"
    "G29" = "5"
    "G30" = "5"
    "G31" = "5"
    "G32" = "5"
    "G33" = "5"
    "G34" = "5"
    "G35" = "5"
    "G36" = "5"
    "G37" = "5"
    "D38" = "0.02417"
    "E38" = "-1.88%"
    "G38" = "5"
    "D39" = "0.05158"
    "E39" = "-2.23%"
    "G39" = "5"
    "D40" = "0.007750"
    "E40" = "-1.20%"
    "G40" = "5"
    "E41" = "-1.11%"
    "G41" = "5"
    "D42" = "0.007625"
    "E42" = "11.63%"
    "G42" = "5"
    "E43" = "72.40%"
    "G43" = "5"
    "D44" = "0.008180"
    "E44" = "5.73%"
    "G44" = "5"
    "D45" = "0.3100"
    "E45" = "1.69%"
    "G45" = "5"
    "D46" = "0.00006583"
    "E46" = "-2.29%"
    "G46" = "5"
    "E47" = "-0.22%"
    "G47" = "5"
    "D48" = "0.2690"
    "E48" = "58.39%"
    "G48" = "5"
    "G49" = "5"
    "E50" = "-0.22%"
    "G50" = "5"
    "E51" = "-0.22%"
    "G51" = "5"
}

foreach ($addr in $cells.Keys) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $cells[$addr]
    $range.Style = "Normal"
}
